# Finished the presentation examples and added screenshot method to base
# The "Search Data" sheet was reworked from the old Incorta example to a
# Selenium example, and the now-unused "Menu" columns (D:E) were dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: swap the Incorta sample data for the Selenium sample data.
$ws.Range("B2").Value = "Selenium"
$ws.Range("C2").Value = "What is Selenium? Introduction to Selenium Automation Testing"

# Columns D:E ("Menu" / "Menu Item") are no longer used - clear them out.
$ws.Range("D1:E2").ClearContents()

# Widen column C to fit the new (longer) text and drop the old autofit flag.
$ws.Columns.Item(3).ColumnWidth = 58.02

# Leave the selection on C2, matching the last-edited cell.
$ws.Range("C2").Select()
